$wb = $excel.ActiveWorkbook

# --- Add the new "Grouping" worksheet as the last sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Grouping"

# --- Populate the grouping analysis table ---------------------------------
$newSheet.Range("A1").Value = "Group No."
$newSheet.Range("B1").Value = "Description"
$newSheet.Range("C1").Value = "Mask"
$newSheet.Range("D1").Value = "Property"

$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "Inactive subscribers"
$newSheet.Range("C2").Value = "df_whizz1.active==0"
$newSheet.Range("D2").Value = "No lesson activity records."

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "Active subscribers taking no assessment"
$newSheet.Range("C3").Value = "(df_whizz1.active==1) & (df_whizz1.assess==0)"
$newSheet.Range("D3").Value = "No performance measure/ academic feedback from Whizz system."

$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "Other active subscribers"
$newSheet.Range("C4").Value = "(df_whizz1.active==1) & (df_whizz1.assess==1)"

# Widen the text columns so the table reads like the original note.
$newSheet.Columns.Item(2).ColumnWidth = 36
$newSheet.Columns.Item(3).ColumnWidth = 39.5546875
$newSheet.Columns.Item(4).ColumnWidth = 56.33203125

# --- Update the previously-active TODO sheet's remembered selection -------
$todo = $wb.Worksheets.Item("TODO")
$todo.Activate()
$todo.Range("A18").Select() | Out-Null

# --- Re-activate "Grouping" so it becomes the active/selected tab ---------
$newSheet.Activate()
$newSheet.Range("D6").Select() | Out-Null
